$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '69.457.63'
Set-TextValue $ws.Range("E2") '  +1.76%  '
Set-TextValue $ws.Range("D3") '3.945.89'
Set-TextValue $ws.Range("E3") '  +0.54%  '
Set-TextValue $ws.Range("D5") '509.04'
Set-TextValue $ws.Range("E5") '  +4.66%  '
Set-TextValue $ws.Range("D6") '147.49'
Set-TextValue $ws.Range("E6") '  -0.21%  '
Set-TextValue $ws.Range("D7") '0.624'
Set-TextValue $ws.Range("E7") '  +0.06%  '
Set-TextValue $ws.Range("D8") '0.999'
Set-TextValue $ws.Range("E8") '  +0.09%  '
Set-TextValue $ws.Range("D9") '0.734'
Set-TextValue $ws.Range("E9") '  +0.11%  '
Set-TextValue $ws.Range("D10") '0.174'
Set-TextValue $ws.Range("E10") '  +4.69%  '
Set-TextValue $ws.Range("D11") '0.0000346'
Set-TextValue $ws.Range("E11") '  -0.91%  '
Set-TextValue $ws.Range("D12") '43.46'
Set-TextValue $ws.Range("E12") '  +0.96%  '
Set-TextValue $ws.Range("D13") '10.49'
Set-TextValue $ws.Range("E13") '  -2.22%  '
Set-TextValue $ws.Range("D14") '4.580.88'
Set-TextValue $ws.Range("E14") '  +0.75%  '
Set-TextValue $ws.Range("D15") '3.981.20'
Set-TextValue $ws.Range("E15") '  +1.66%  '
Set-TextValue $ws.Range("D16") '14.20'
Set-TextValue $ws.Range("E16") '  -1.57%  '
Set-TextValue $ws.Range("E17") '  -0.35%  '
Set-TextValue $ws.Range("E18") '  +7.71%  '
Set-TextValue $ws.Range("D19") '19.93'
Set-TextValue $ws.Range("E19") '  -0.02%  '
Set-TextValue $ws.Range("D20") '69.553.64'
Set-TextValue $ws.Range("E20") '  +1.77%  '
Set-TextValue $ws.Range("D21") '434.67'
Set-TextValue $ws.Range("E21") '  -1.77%  '
Set-TextValue $ws.Range("D22") '3.43'
Set-TextValue $ws.Range("E22") '  -1.34%  '
Set-TextValue $ws.Range("D23") '14.60'
Set-TextValue $ws.Range("E23") '  -3.99%  '
Set-TextValue $ws.Range("D24") '88.74'
Set-TextValue $ws.Range("E24") '  +0.58%  '
Set-TextValue $ws.Range("D25") '11.77'
Set-TextValue $ws.Range("E25") '  +4.22%  '
Set-TextValue $ws.Range("D26") '3.88'
Set-TextValue $ws.Range("E26") '  +7.27%  '
Set-TextValue $ws.Range("D27") '11.15'
Set-TextValue $ws.Range("E27") '  -2.86%  '
Set-TextValue $ws.Range("D28") '37.06'
Set-TextValue $ws.Range("E28") '  -4.35%  '
Set-TextValue $ws.Range("E29") '  -0.56%  '
Set-TextValue $ws.Range("D30") '708.59'
Set-TextValue $ws.Range("E30") '  -1.36%  '
Set-TextValue $ws.Range("D31") '13.36'
Set-TextValue $ws.Range("E31") '  -3.01%  '
Set-TextValue $ws.Range("D32") '0.128'
Set-TextValue $ws.Range("E32") '  -1.86%  '
Set-TextValue $ws.Range("D33") '2.89'
Set-TextValue $ws.Range("E33") '  -0.72%  '
Set-TextValue $ws.Range("D34") '68.22'
Set-TextValue $ws.Range("E34") '  +12.00%  '
Set-TextValue $ws.Range("D35") '0.443'
Set-TextValue $ws.Range("E35") '  +8.39%  '
Set-TextValue $ws.Range("D36") '0.0₃0878'
Set-TextValue $ws.Range("E36") '  -1.34%  '
Set-TextValue $ws.Range("D37") '6.01'
Set-TextValue $ws.Range("E37") '  -6.71%  '
Set-TextValue $ws.Range("D38") '40.81'
Set-TextValue $ws.Range("E38") '  -4.09%  '
Set-TextValue $ws.Range("D39") '0.149'
Set-TextValue $ws.Range("E39") '  -0.78%  '
Set-TextValue $ws.Range("E40") '  -0.06%  '
Set-TextValue $ws.Range("E41") '  +0.01%  '
Set-TextValue $ws.Range("D42") '0.0490'
Set-TextValue $ws.Range("E42") '  +1.68%  '
Set-TextValue $ws.Range("D43") '2.87'
Set-TextValue $ws.Range("E43") '  -4.42%  '
Set-TextValue $ws.Range("D44") '3.15'
Set-TextValue $ws.Range("E44") '  +8.40%  '
Set-TextValue $ws.Range("D45") '3.06'
Set-TextValue $ws.Range("E45") '  -6.47%  '
Set-TextValue $ws.Range("D46") '0.144'
Set-TextValue $ws.Range("E46") '  +1.21%  '
Set-TextValue $ws.Range("B47") 'BabyDogeCoin'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D47") '0.0₆0366'
Set-TextValue $ws.Range("E47") '  +1.90%  '
Set-TextValue $ws.Range("B48") 'ApeXProtocol'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D48") '3.35'
Set-TextValue $ws.Range("E48") '  +1.04%  '
Set-TextValue $ws.Range("D49") '3.00'
Set-TextValue $ws.Range("E49") '  +5.48%  '
Set-TextValue $ws.Range("D50") '3.39'
Set-TextValue $ws.Range("E50") '  -1.19%  '
Set-TextValue $ws.Range("D51") '2.11'
Set-TextValue $ws.Range("E51") '  -1.40%  '
